# Applies the edits described by the commit diff:
#  - Produits: E4 (stock) 12 -> 8
#  - Ventes: append two new sale rows (5 and 6)
#  - Clients: C3 (Dette) 50000 -> 49000 ; D3 (Total Achats) 50000 -> 150000

$wb = $excel.ActiveWorkbook

# --- Sheet "Produits": update stock for product in row 4 ---
$wsProduits = $wb.Worksheets.Item("Produits")
$wsProduits.Range("E4").Value = 8

# --- Sheet "Ventes": append two new sales rows ---
$wsVentes = $wb.Worksheets.Item("Ventes")
$wsVentes.Range("A5").Value = "VNT-20251223011259-1"
$wsVentes.Range("B5").Value = 100000
$wsVentes.Range("C5").Value = "cash"
$wsVentes.Range("D5").Value = "2025-12-23 00:12:59"

$wsVentes.Range("A6").Value = "VNT-20251223011405-1"
$wsVentes.Range("B6").Value = 100000
$wsVentes.Range("C6").Value = "cash"
$wsVentes.Range("D6").Value = "2025-12-23 00:14:05"

# --- Sheet "Clients": update debt and total purchases for client in row 3 ---
$wsClients = $wb.Worksheets.Item("Clients")
$wsClients.Range("C3").Value = 49000
$wsClients.Range("D3").Value = 150000
